# The single worksheet's name encodes the date the approved-systems list was
# last updated. This edit bumps that date from 02-12-2025 to 05-12-2025 by
# renaming the sheet accordingly.
#
# Excel keeps the workbook-level defined name
# "Telemedicinsk_hjemmemonitoriering" (which refers to
# '<sheet name>'!$A$1:$D$10) automatically in sync with the sheet's name,
# so renaming the sheet alone is sufficient to update the defined name too.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$newName = "Opdateret d. 05-12-2025"

if ($ws.Name -ne $newName) {
    $ws.Name = $newName
}
